$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowByName($name) {
    $cell = $ws.Cells.Find($name)
    return $cell.Row()
}

# 1. Update the report title/date in the header cell (A1).
$ws.Range("A1").Value = "Item Name 1930 2/9/25"

# 2. Update existing prices that changed.
$ws.Range("D" + (Get-RowByName("Brussel Sprouts"))).Value = 29.75
$ws.Range("D" + (Get-RowByName("Grapes - Red (Seedless)"))).Value = 35
$ws.Range("D" + (Get-RowByName("Leeks"))).Value = 44.25
$ws.Range("D" + (Get-RowByName("Potato - Russet"))).Value = 23.75
$ws.Range("D" + (Get-RowByName("Strawberries Fresh"))).Value = 34.75
$ws.Range("D" + (Get-RowByName("Zucchini Fancy Medium"))).Value = 28.5

# 3. Update case-size labels for the green bell pepper rows ("case" -> "1 cs").
$ws.Range("C" + (Get-RowByName("Pepper - Green Bell (Extra Large)"))).Value = "1 cs"
$pepperLargeRow = Get-RowByName("Pepper - Green Bell (Large)")
$ws.Range("C" + $pepperLargeRow).Value = "1 cs"
$ws.Range("D" + $pepperLargeRow).Value = 24.75

# 4. Insert a new "Melon - Honeydew" row right after "Melon - Cantaloupe",
#    copying the row-above formatting so styles match the rest of the table.
$cantaloupeRow = Get-RowByName("Melon - Cantaloupe")
$newRow = $cantaloupeRow + 1
$ws.Rows.Item($newRow).Insert()
$ws.Range("A" + $cantaloupeRow + ":D" + $cantaloupeRow).Copy()
$ws.Range("A" + $newRow + ":D" + $newRow).PasteSpecial(-4122)

$ws.Range("A" + $newRow).Value = "Melon - Honeydew"
$ws.Range("B" + $newRow).Value = 1425
$ws.Range("C" + $newRow).Value = "1 cs"
$ws.Range("D" + $newRow).Value = 33.75

$excel.CutCopyMode = 0
